$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 96, pushing the existing rows 96-154
# (which already carry all the correct data for what will become rows
# 98-156) down to rows 98-156.
$ws.Rows("96:97").Insert()

# Populate the two brand-new rows (96 and 97) with the new weekly price
# records for "Femacal de La Calera - Alcachofa".

# Row 96
$ws.Range("A96").Value = 3
$ws.Range("B96").Value = "Femacal de La Calera"
$ws.Range("C96").Value = "Coquimbo"
$ws.Range("D96").Value = 44438
$ws.Range("E96").Value = 5
$ws.Range("F96").Value = 100112013
$ws.Range("G96").Value = "Alcachofa"
$ws.Range("H96").Value = "Argentina(o)"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 140
$ws.Range("K96").Value = 11000
$ws.Range("L96").Value = 11500
$ws.Range("M96").Value = 11286
$ws.Range("N96").Value = "$/caja 50 unidades"
$ws.Range("O96").Value = "Provincia de Limarí"
$ws.Range("P96").Value = 226
$ws.Range("Q96").Value = 50
$ws.Range("R96").Value = "Hortaliza"

# Row 97
$ws.Range("A97").Value = 3
$ws.Range("B97").Value = "Femacal de La Calera"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44438
$ws.Range("E97").Value = 5
$ws.Range("F97").Value = 100112013
$ws.Range("G97").Value = "Alcachofa"
$ws.Range("H97").Value = "Española"
$ws.Range("I97").Value = "Extra"
$ws.Range("J97").Value = 130
$ws.Range("K97").Value = 12000
$ws.Range("L97").Value = 12500
$ws.Range("M97").Value = 12231
$ws.Range("N97").Value = "$/caja 30 unidades"
$ws.Range("O97").Value = "Provincia de Limarí"
$ws.Range("P97").Value = 408
$ws.Range("Q97").Value = 30
$ws.Range("R97").Value = "Hortaliza"

# Make sure the date cells keep the same date/time number format used by
# the rest of column D.
$ws.Range("D96:D97").NumberFormat = $ws.Range("D98").NumberFormat
